# Switch example group numbers:
# Replace the example team roster (Alice/Bob/Claire/David/Elaine) with the
# real group's members (Veselin/Rawda/Hannah/Mirit/Bogdana), and add the
# sixth member (Martin) as a new row, copying the rating pattern used by
# the first renamed row.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

# Rename the five existing team members.
$ws.Range("B12").Value = "Veselin"
$ws.Range("B13").Value = "Rawda"
$ws.Range("B14").Value = "Hannah"
$ws.Range("B15").Value = "Mirit"
$ws.Range("B16").Value = "Bogdana"

# Add the new sixth team member as row 17.
$ws.Range("B17").Value = "Martin"
$ws.Range("C17").Value = 4
$ws.Range("E17").Value = 1

# All six name cells get the new explicit black font color.
$ws.Range("B12:B17").Font.Color = 0

# Move the active selection to reflect where the author left off editing.
[void]$ws.Range("F14").Select()
